$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Make room for the new rows.
#    * 7 new rows before the old row 5 (p21 q4/q5 + p22 q1-5)
#    * 2 new rows before the old p24-q5 row (now row 14) for p24 q3/q4
# ------------------------------------------------------------------
$ws.Rows("5:11").Insert(-4121)    # xlShiftDown
$ws.Rows("14:15").Insert(-4121)   # xlShiftDown

# ------------------------------------------------------------------
# 2) Fill in the "quest file" column for the existing p21 rows that
#    were missing it.
# ------------------------------------------------------------------
$ws.Range("G2").Value = "1.html"
$ws.Range("G3").Value = "2.html"

# ------------------------------------------------------------------
# 3) p21 q4 and q5 (rows 5 and 6) - new rows.
# ------------------------------------------------------------------
$ws.Range("A5").Value = "econ"
$ws.Range("B5").Value = 2025
$ws.Range("C5").Value = "AS"
$ws.Range("D5").Value = "qp-202505-economics-p21"
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = "subjects/econ/2025/AS/qp-202505-economics-p21"
$ws.Range("G5").Value = "4.html"

$ws.Range("A6").Value = "econ"
$ws.Range("B6").Value = 2025
$ws.Range("C6").Value = "AS"
$ws.Range("D6").Value = "qp-202505-economics-p21"
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = "subjects/econ/2025/AS/qp-202505-economics-p21"
$ws.Range("G6").Value = "5.html"

# ------------------------------------------------------------------
# 4) p22 q1-q5 (rows 7-11) - brand new question paper block.
# ------------------------------------------------------------------
for ($q = 1; $q -le 5; $q++) {
    $r = 6 + $q
    $ws.Range("A$r").Value = "econ"
    $ws.Range("B$r").Value = 2025
    $ws.Range("C$r").Value = "AS"
    $ws.Range("D$r").Value = "qp-202505-economics-p22"
    $ws.Range("E$r").Value = $q
    $ws.Range("F$r").Value = "subjects/econ/2025/AS/qp-202505-economics-p22"
    $ws.Range("G$r").Value = "$q.html"
}

# ------------------------------------------------------------------
# 5) p24 q1 and q2 (rows 12 and 13) already hold the subject/year/
#    qp info (shifted down from the original sheet) - just add the
#    "quest file" value, matching the styling used in the rest of
#    the p24 block (copy format from row 16, the original p24 q5
#    row which already carries that styling).
# ------------------------------------------------------------------
$ws.Range("G16").Copy()
$ws.Range("G12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G12").Value = "1.html"

$ws.Range("G16").Copy()
$ws.Range("G13").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G13").Value = "2.html"

# ------------------------------------------------------------------
# 6) p24 q3 and q4 (rows 14 and 15) - brand new rows. Only the
#    "quest file" (G) and "new" (H) columns carry the special
#    styling used by the rest of the p24 block; A-F stay default.
# ------------------------------------------------------------------
$ws.Range("A14").Value = "econ"
$ws.Range("B14").Value = 2025
$ws.Range("C14").Value = "AS"
$ws.Range("D14").Value = "qp-202505-economics-p24"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = "subjects/econ/2025/AS/qp-202505-economics-p24"

$ws.Range("G16").Copy()
$ws.Range("G14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G14").Value = "3.html"

$ws.Range("H16").Copy()
$ws.Range("H14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H14").ClearContents()

$ws.Range("A15").Value = "econ"
$ws.Range("B15").Value = 2025
$ws.Range("C15").Value = "AS"
$ws.Range("D15").Value = "qp-202505-economics-p24"
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = "subjects/econ/2025/AS/qp-202505-economics-p24"

$ws.Range("G16").Copy()
$ws.Range("G15").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G15").Value = "4.html"

$ws.Range("H16").Copy()
$ws.Range("H15").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H15").ClearContents()

# ------------------------------------------------------------------
# 7) p24 q5 (row 16) keeps its original values - nothing to change.
# ------------------------------------------------------------------

# ------------------------------------------------------------------
# 8) p24 q6 (row 17) - new row, reuse row 16's styling for columns
#    A-G (column H / "new" stays completely empty, no cell at all).
# ------------------------------------------------------------------
$ws.Range("A16:G16").Copy()
$ws.Range("A17:G17").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A17").Value = "econ"
$ws.Range("B17").Value = 2025
$ws.Range("C17").Value = "AS"
$ws.Range("D17").Value = "qp-202505-economics-p24"
$ws.Range("E17").Value = 6
$ws.Range("F17").Value = "subjects/econ/2025/AS/qp-202505-economics-p24"
$ws.Range("G17").Value = "6.html"

# ------------------------------------------------------------------
# 9) Selection / view state, matching the final saved workbook.
# ------------------------------------------------------------------
$ws.Range("E12").Select()
